# Junction_Flooding_135.xlsx edit:
#  - reduce "custom accuracy" of row 5 data to 2 decimal places
#  - remove the last data row (row 6) -> "데이터 1000개" (trim dataset)
#  - three data columns (B, O, V) get one unit narrower (width 8 -> 7)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Narrow columns B (2), O (15) and V (22) from width 8 to width 7 ---
# (COM ColumnWidth = XML width - 0.83, so target 7 -> 6.17)
$ws.Columns.Item(2).ColumnWidth = 6.17
$ws.Columns.Item(15).ColumnWidth = 6.17
$ws.Columns.Item(22).ColumnWidth = 6.17

# --- Round row 5 (B5:AH5) values to 2 decimal places ---
$row5Values = @(
    15.37, 11.21, 1.07, 33.13, 27.4, 12.1, 45.76, 18.62, 8.2, 12.28,
    13.39, 13.99, 3.86, 12.03, 17.07, 10.2, 0.82, 0.68, 175.79, 33.62,
    11.1, 22.51, 12.02, 1.5, 22.33, 9.81, 8.75, 10.28, 14.0, 0.56,
    41.37, 6.22, 13.88
)

for ($i = 0; $i -lt $row5Values.Length; $i++) {
    # column 2 = B .. column 34 = AH
    $ws.Cells.Item(5, 2 + $i).Value = $row5Values[$i]
}

# --- Remove the last data row (row 6), shifting dimension to A1:AH5 ---
$ws.Rows.Item(6).Delete()
